$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.934.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.882.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.020"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.71%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.018"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3900"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.82"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07908"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.887.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.901"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.065"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.021"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06757"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.019"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.947.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.442"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.351"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.115.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.23%  "

$ws.Range("E28").Value = "  -2.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.051"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.417"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09450"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9539"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.665"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.295"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.343"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06076"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02224"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.209"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.057"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5902"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1870"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.50%  "

$ws.Range("E44").Value = "  +1.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5610"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.395"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.901"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06890"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.45%  "

$ws.Range("E51").Value = "  -1.22%  "
